$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new CRUD test data cells (rows 14-33)
$ws.Range("A14").Value = "Pruebas CRUD Clientes"
$ws.Range("B14").Value = "Cedula"
$ws.Range("C14").Value = "Nombre"
$ws.Range("D14").Value = "Apellido"
$ws.Range("E14").Value = "Genero"
$ws.Range("F14").Value = "Categoria"
$ws.Range("G14").Value = "Resultado"
$ws.Range("A16").Value = "Insertar"
$ws.Range("A18").Value = "Cedula solo letras"
$ws.Range("B18").Value = "prueba"
$ws.Range("C18").Value = "prueba"
$ws.Range("D18").Value = "prueba"
$ws.Range("E18").Value = "masculino"
$ws.Range("F18").Value = "premium"
$ws.Range("G18").Value = "no se registró"
$ws.Range("A19").Value = "Nombre con numeros"
$ws.Range("B19").Value = "prueba1"
$ws.Range("C19").Value = "prueba1"
$ws.Range("D19").Value = "prueba"
$ws.Range("E19").Value = "masculino"
$ws.Range("F19").Value = "normal"
$ws.Range("G19").Value = "no se registró"
$ws.Range("A20").Value = "Cedula con numeros"
$ws.Range("B20").Value = "prueba1"
$ws.Range("C20").Value = "prueba"
$ws.Range("D20").Value = "prueba1"
$ws.Range("E20").Value = "femenino"
$ws.Range("F20").Value = "premium"
$ws.Range("G20").Value = "no se registró"
$ws.Range("A21").Value = "Campo obligatorio faltante"
$ws.Range("B21").Value = "prueba1"
$ws.Range("C21").Value = "prueba"
$ws.Range("D21").Value = "prueba1"
$ws.Range("E21").Value = "."
$ws.Range("F21").Value = "premium"
$ws.Range("G21").Value = "no se registró"
$ws.Range("A22").Value = "Registro correcto"
$ws.Range("B22").Value = "prueba1"
$ws.Range("C22").Value = "nombre"
$ws.Range("D22").Value = "prueba"
$ws.Range("E22").Value = "masculino"
$ws.Range("F22").Value = "premium"
$ws.Range("G22").Value = "se registró"
$ws.Range("A24").Value = "Mostrar"
$ws.Range("C24").Value = "Resultado"
$ws.Range("A26").Value = "Cedula registrada"
$ws.Range("B26").Value = "cedula1"
$ws.Range("C26").Value = "se muestra"
$ws.Range("A27").Value = "Cedula registrada"
$ws.Range("B27").Value = "cedula2"
$ws.Range("C27").Value = "se muestra"
$ws.Range("A28").Value = "Cedula registrada"
$ws.Range("B28").Value = "cedula3"
$ws.Range("C28").Value = "se muestra"
$ws.Range("A29").Value = "Cedula no registrada"
$ws.Range("B29").Value = "cedula4"
$ws.Range("C29").Value = "no se muestra"
$ws.Range("A31").Value = "Editar/Eliminar"
$ws.Range("F31").Value = "ResultadoEditar"
$ws.Range("G31").Value = "ResultadoEliminar"
$ws.Range("A33").Value = "Cambiado"
$ws.Range("B33").Value = "cedulacamb1"
$ws.Range("C33").Value = "cambiado"
$ws.Range("D33").Value = "cambiado"
$ws.Range("E33").Value = "femenino"
$ws.Range("F33").Value = "actualizó"
$ws.Range("G33").Value = "se eliminó"

# Column widths (character units); engine rounds to nearest 1/6 + 5/6 offset internally
$ws.Columns.Item(1).ColumnWidth = 24.166666666666668   # -> stored 25.0
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666   # -> stored 14.0
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666   # -> stored 14.0
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334   # -> stored ~15.67 (closest to 15.7109375)
$ws.Columns.Item(6).ColumnWidth = 14.166666666666666   # -> stored 15.0
$ws.Columns.Item(7).ColumnWidth = 16.5                 # -> stored ~17.33 (closest to 17.28515625)

# Selection / view state
$ws.Range("C22").Select()

# Page setup orientation
$ws.PageSetup.Orientation = 1
